# Scheduled-runner price/profit refresh across the Leve profit sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the derived
# LeveProfit(NQ/HQ) columns (H,I,J,K,L,M,N) for the rows whose market
# data changed in this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 29251.295
$ws.Range("J128").Value = 29251.295
$ws.Range("L128").Value = 29251.295
$ws.Range("N128").Value = -39211.295

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4458.172
$ws.Range("I32").Value = 4000.7292
$ws.Range("J32").Value = 5830.5
$ws.Range("K32").Value = 4000.7292
$ws.Range("L32").Value = 5830.5
$ws.Range("M32").Value = -3713.7292
$ws.Range("N32").Value = -6404.5

$ws.Range("H132").Value = 2360.6
$ws.Range("I132").Value = 2121.75
$ws.Range("J132").Value = 3316
$ws.Range("K132").Value = 6365.25
$ws.Range("L132").Value = 9948
$ws.Range("M132").Value = -3835.25
$ws.Range("N132").Value = -15008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1634.2903
$ws.Range("I20").Value = 1423.0952
$ws.Range("K20").Value = 1423.0952
$ws.Range("M20").Value = -1176.0952

$ws.Range("H86").Value = 1960.2667
$ws.Range("I86").Value = 1745.8182
$ws.Range("J86").Value = 2550
$ws.Range("K86").Value = 1745.8182
$ws.Range("L86").Value = 2550
$ws.Range("M86").Value = -622.8181999999999
$ws.Range("N86").Value = -4796

$ws.Range("H89").Value = 1960.2667
$ws.Range("I89").Value = 1745.8182
$ws.Range("J89").Value = 2550
$ws.Range("K89").Value = 8729.091
$ws.Range("L89").Value = 12750
$ws.Range("M89").Value = -3113.091
$ws.Range("N89").Value = -23982

$ws.Range("H134").Value = 2658.303
$ws.Range("I134").Value = 2274.9565
$ws.Range("J134").Value = 3540
$ws.Range("K134").Value = 6824.869499999999
$ws.Range("L134").Value = 10620
$ws.Range("M134").Value = -4289.869499999999
$ws.Range("N134").Value = -15690

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 79800
$ws.Range("J20").Value = 79800
$ws.Range("L20").Value = 79800
$ws.Range("N20").Value = -80272

$ws.Range("H30").Value = 79800
$ws.Range("J30").Value = 79800
$ws.Range("L30").Value = 79800
$ws.Range("N30").Value = -79982

$ws.Range("H31").Value = 3207.9058
$ws.Range("I31").Value = 3200.8333
$ws.Range("J31").Value = 3211.543
$ws.Range("K31").Value = 3200.8333
$ws.Range("L31").Value = 3211.543
$ws.Range("M31").Value = -2905.8333
$ws.Range("N31").Value = -3801.543

$ws.Range("H34").Value = 3207.9058
$ws.Range("I34").Value = 3200.8333
$ws.Range("J34").Value = 3211.543
$ws.Range("K34").Value = 3200.8333
$ws.Range("L34").Value = 3211.543
$ws.Range("M34").Value = -2998.8333
$ws.Range("N34").Value = -3615.543

$ws.Range("H58").Value = 2213.7778
$ws.Range("I58").Value = 2184.8
$ws.Range("J58").Value = 2250
$ws.Range("K58").Value = 2184.8
$ws.Range("L58").Value = 2250
$ws.Range("M58").Value = -1981.8
$ws.Range("N58").Value = -2656

$ws.Range("H107").Value = 689.6
$ws.Range("I107").Value = 475.25
$ws.Range("J107").Value = 832.5
$ws.Range("K107").Value = 475.25
$ws.Range("L107").Value = 832.5
$ws.Range("M107").Value = 1444.75
$ws.Range("N107").Value = -4672.5

$ws.Range("H128").Value = 79800
$ws.Range("J128").Value = 79800
$ws.Range("L128").Value = 79800
$ws.Range("N128").Value = -89760

$ws.Range("H134").Value = 6871.65
$ws.Range("I134").Value = 8218.75
$ws.Range("J134").Value = 4851
$ws.Range("K134").Value = 24656.25
$ws.Range("L134").Value = 14553
$ws.Range("M134").Value = -22121.25
$ws.Range("N134").Value = -19623

$ws.Range("H136").Value = 2213.7778
$ws.Range("I136").Value = 2184.8
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 6554.400000000001
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -4004.400000000001
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 6500
$ws.Range("J55").Value = 6500
$ws.Range("L55").Value = 19500
$ws.Range("N55").Value = -19854

$ws.Range("H64").Value = 335170.34
$ws.Range("I64").Value = 2012
$ws.Range("J64").Value = 501749.5
$ws.Range("K64").Value = 6036
$ws.Range("L64").Value = 1505248.5
$ws.Range("M64").Value = -5766
$ws.Range("N64").Value = -1505788.5

$ws.Range("H67").Value = 335170.34
$ws.Range("I67").Value = 2012
$ws.Range("J67").Value = 501749.5
$ws.Range("K67").Value = 6036
$ws.Range("L67").Value = 1505248.5
$ws.Range("M67").Value = -5100
$ws.Range("N67").Value = -1507120.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5416.3135
$ws.Range("I70").Value = 5557.8037
$ws.Range("K70").Value = 5557.8037
$ws.Range("M70").Value = -5287.8037

$ws.Range("H73").Value = 5416.3135
$ws.Range("I73").Value = 5557.8037
$ws.Range("K73").Value = 5557.8037
$ws.Range("M73").Value = -4621.8037

$ws.Range("H97").Value = 556.1539
$ws.Range("I97").Value = 556.1539
$ws.Range("K97").Value = 556.1539
$ws.Range("M97").Value = -60.15390000000002

$ws.Range("H102").Value = 1116.2
$ws.Range("I102").Value = 695.25
$ws.Range("K102").Value = 695.25
$ws.Range("M102").Value = 926.75

$ws.Range("H132").Value = 5858.706
$ws.Range("I132").Value = 9166.333000000001
$ws.Range("J132").Value = 4054.5454
$ws.Range("K132").Value = 27498.999
$ws.Range("L132").Value = 12163.6362
$ws.Range("M132").Value = -24968.999
$ws.Range("N132").Value = -17223.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 64916.938
$ws.Range("I7").Value = 93052.82000000001
$ws.Range("J7").Value = 3018
$ws.Range("K7").Value = 93052.82000000001
$ws.Range("L7").Value = 3018
$ws.Range("M7").Value = -92940.82000000001
$ws.Range("N7").Value = -3242

$ws.Range("H16").Value = 4597.4
$ws.Range("J16").Value = 5493.5
$ws.Range("L16").Value = 5493.5
$ws.Range("N16").Value = -5833.5

$ws.Range("H46").Value = 20835154
$ws.Range("I46").Value = 37038260
$ws.Range("J46").Value = 2589.1428
$ws.Range("K46").Value = 37038260
$ws.Range("L46").Value = 2589.1428
$ws.Range("M46").Value = -37038072
$ws.Range("N46").Value = -2965.1428

$ws.Range("H59").Value = 27392.5
$ws.Range("J59").Value = 27392.5
$ws.Range("L59").Value = 27392.5
$ws.Range("N59").Value = -28700.5

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H126").Value = 64916.938
$ws.Range("I126").Value = 93052.82000000001
$ws.Range("J126").Value = 3018
$ws.Range("K126").Value = 279158.46
$ws.Range("L126").Value = 9054
$ws.Range("M126").Value = -276688.46
$ws.Range("N126").Value = -13994

$ws.Range("H132").Value = 20640486
$ws.Range("I132").Value = 25496012
$ws.Range("K132").Value = 76488036
$ws.Range("M132").Value = -76485506

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2932.9524
$ws.Range("I132").Value = 2818.875
$ws.Range("J132").Value = 3298
$ws.Range("K132").Value = 8456.625
$ws.Range("L132").Value = 9894
$ws.Range("M132").Value = -5926.625
$ws.Range("N132").Value = -14954
